# Update the scraped_date column (D) on every worksheet from the old
# timestamp to the new one, leaving every other cell untouched.

$wb = $excel.ActiveWorkbook

$oldDate = "2026-02-21 21:19:38"
$newDate = "2026-02-21 21:23:34"

foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Rows.Count

    # Row 1 is the header ("scraped_date"); data starts on row 2.
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 4)
        if ($cell.Value2 -eq $oldDate) {
            $cell.Value2 = $newDate
        }
    }
}
